# Automatische test-sync: 2025-07-31 21:39:50
# Adds Testmail #8 as a new row (row 10) to the "Logs" sheet, extends the
# conditional-formatting ranges to cover it, and refreshes the category
# pivot on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "Logs" sheet: append the new row of mail-log data
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A10").Value = "Kun je nagaan of dit nog leverbaar is?"
$logs.Range("B10").Value = "mailmind.test@zohomail.eu"
$logs.Range("C10").Value = "Testmail #8: Kun je nagaan of dit nog leverbaar is?"
$logs.Range("D10").Value = "Productinformatie"
$logs.Range("E10").Value = "Beste klant,`r`nDank voor uw e-mail. Om u beter van dienst te kunnen zijn, heb ik meer specifieke informatie nodig over het product waar u naar informeert. Kunt u ons het productnummer, de naam van het product of enige details geven zodat we het kunnen controleren in ons systeem? `r`nMet vriendelijke groet,`r`n[Bedrijfsnaam]"
$logs.Range("F10").Value = "2025-07-31 21:39:04"
$logs.Range("G10").Value = "Ja"
$logs.Range("H10").Value = "Nee"
$logs.Range("I10").Value = "Ja"
$logs.Range("J10").Value = "Nee"

# Multi-line content in E10 triggers an auto row-height bump; restore the
# default (no explicit row height, matching the other data rows).
$logs.Rows.Item(10).AutoFit()

# ---------------------------------------------------------------------
# 2. Extend the conditional-formatting ranges (D/G/H/I/J) from row 9 to
#    row 10 so the newly added row is covered too.
# ---------------------------------------------------------------------
$logs.Range("D2:D10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D10"))
$logs.Range("G2:G10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G10"))
$logs.Range("H2:H10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H10"))
$logs.Range("I2:I10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I10"))
$logs.Range("J2:J10").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J10"))

# ---------------------------------------------------------------------
# 3. "Dashboard" sheet: the category counts behind the chart move --
#    "Productinformatie" now leads with 3 hits, "Intern verzoek / Actie
#    voor medewerker" drops to the row below with its original count.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Productinformatie"
$dash.Range("B3").Value = 3

$dash.Range("A4").Value = "Intern verzoek / Actie voor medewerker"
$dash.Range("B4").Value = 2
